$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "codeforiati:group-code" (col C) and "codeforiati:group-name" (col D)
# columns need to swap places so that group-name comes before group-code:
# column C becomes the name column and column D becomes the code column.
# Swap the two ranges using a scratch column (F, unused in this sheet) so we
# don't disturb column-level properties via a whole-column cut/insert.
$ws.Range("D1:D94").Copy()
$ws.Range("F1").PasteSpecial()
$ws.Range("C1:C94").Copy()
$ws.Range("D1").PasteSpecial()
$ws.Range("F1:F94").Copy()
$ws.Range("C1").PasteSpecial()
$ws.Range("F1:F94").ClearContents()
